$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="43.075.17"; ForceText=$false},
    @{Cell="D3"; Value="2.306.27"; ForceText=$false},
    @{Cell="E3"; Value="  -0.32%  "; ForceText=$false},
    @{Cell="D5"; Value="300.97"; ForceText=$true},
    @{Cell="E5"; Value="  -0.40%  "; ForceText=$false},
    @{Cell="D6"; Value="98.60"; ForceText=$true},
    @{Cell="E6"; Value="  -2.91%  "; ForceText=$false},
    @{Cell="E7"; Value="  +3.81%  "; ForceText=$false},
    @{Cell="E8"; Value="  +0.05%  "; ForceText=$false},
    @{Cell="D9"; Value="0.521"; ForceText=$true},
    @{Cell="E9"; Value="  +1.16%  "; ForceText=$false},
    @{Cell="D10"; Value="35.68"; ForceText=$true},
    @{Cell="E10"; Value="  -1.73%  "; ForceText=$false},
    @{Cell="E11"; Value="  -0.18%  "; ForceText=$false},
    @{Cell="E12"; Value="  -0.46%  "; ForceText=$false},
    @{Cell="D13"; Value="17.88"; ForceText=$true},
    @{Cell="E13"; Value="  -0.02%  "; ForceText=$false},
    @{Cell="E14"; Value="  +0.30%  "; ForceText=$false},
    @{Cell="D15"; Value="2.663.57"; ForceText=$false},
    @{Cell="D16"; Value="2.345.74"; ForceText=$false},
    @{Cell="E16"; Value="  +1.46%  "; ForceText=$false},
    @{Cell="E17"; Value="  -2.29%  "; ForceText=$false},
    @{Cell="D18"; Value="42.981.85"; ForceText=$false},
    @{Cell="D19"; Value="13.51"; ForceText=$true},
    @{Cell="E19"; Value="  +7.70%  "; ForceText=$false},
    @{Cell="E20"; Value="  +0.59%  "; ForceText=$false},
    @{Cell="E21"; Value="  -1.37%  "; ForceText=$false},
    @{Cell="D22"; Value="68.38"; ForceText=$true},
    @{Cell="E22"; Value="  +0.70%  "; ForceText=$false},
    @{Cell="D23"; Value="239.13"; ForceText=$true},
    @{Cell="E23"; Value="  +1.07%  "; ForceText=$false},
    @{Cell="E24"; Value="  -2.08%  "; ForceText=$false},
    @{Cell="D25"; Value="0.998"; ForceText=$true},
    @{Cell="E25"; Value="  -0.18%  "; ForceText=$false},
    @{Cell="D26"; Value="2.43"; ForceText=$true},
    @{Cell="E26"; Value="  -0.91%  "; ForceText=$false},
    @{Cell="E27"; Value="  +0.25%  "; ForceText=$false},
    @{Cell="D28"; Value="167.68"; ForceText=$true},
    @{Cell="E28"; Value="  -0.69%  "; ForceText=$false},
    @{Cell="E29"; Value="  -0.20%  "; ForceText=$false},
    @{Cell="E30"; Value="  -13.38%  "; ForceText=$false},
    @{Cell="D31"; Value="33.39"; ForceText=$true},
    @{Cell="E31"; Value="  -4.02%  "; ForceText=$false},
    @{Cell="D32"; Value="5.24"; ForceText=$true},
    @{Cell="E32"; Value="  +4.38%  "; ForceText=$false},
    @{Cell="E33"; Value="  -0.05%  "; ForceText=$false},
    @{Cell="D34"; Value="4.82"; ForceText=$true},
    @{Cell="E34"; Value="  +3.13%  "; ForceText=$false},
    @{Cell="D35"; Value="18.17"; ForceText=$true},
    @{Cell="E35"; Value="  +4.10%  "; ForceText=$false},
    @{Cell="E36"; Value="  -0.41%  "; ForceText=$false},
    @{Cell="E37"; Value="  -0.51%  "; ForceText=$false},
    @{Cell="E38"; Value="  -0.69%  "; ForceText=$false},
    @{Cell="E39"; Value="  +0.20%  "; ForceText=$false},
    @{Cell="E40"; Value="  +2.14%  "; ForceText=$false},
    @{Cell="E41"; Value="  -3.25%  "; ForceText=$false},
    @{Cell="D42"; Value="2.005.17"; ForceText=$false},
    @{Cell="E42"; Value="  +0.96%  "; ForceText=$false},
    @{Cell="E43"; Value="  -0.25%  "; ForceText=$false},
    @{Cell="E44"; Value="  -2.42%  "; ForceText=$false},
    @{Cell="D45"; Value="10.08"; ForceText=$true},
    @{Cell="E45"; Value="  -1.43%  "; ForceText=$false},
    @{Cell="D46"; Value="17.44"; ForceText=$true},
    @{Cell="E46"; Value="  -1.01%  "; ForceText=$false},
    @{Cell="E47"; Value="  -3.15%  "; ForceText=$false},
    @{Cell="D48"; Value="54.50"; ForceText=$true},
    @{Cell="E48"; Value="  -2.94%  "; ForceText=$false},
    @{Cell="D49"; Value="2.529.64"; ForceText=$false},
    @{Cell="E49"; Value="  -0.39%  "; ForceText=$false},
    @{Cell="D50"; Value="73.82"; ForceText=$true},
    @{Cell="E50"; Value="  +5.26%  "; ForceText=$false},
    @{Cell="E51"; Value="  +0.20%  "; ForceText=$false}
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    if ($chg.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $chg.Value
}
